$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (volume number and report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Cells changing from a numeric/placeholder style to a text placeholder style (copy style+value from a template cell holding the same placeholder style) ---
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C15").Copy($ws.Range("F15"))
$ws.Range("C15").Copy($ws.Range("C28"))
$ws.Range("C15").Copy($ws.Range("G31"))
$ws.Range("E14").Copy($ws.Range("H31"))

# --- Cells changing from a text placeholder into a plain number (copy style template, then set value) ---
$ws.Range("G15").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K15").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = -12.903225806451
$ws.Range("L16").Value = -12.903225806451
$ws.Range("M16").Value = -47.058823529411
$ws.Range("N16").Value = -81.879194630872
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 45.833333333333
$ws.Range("I17").Value = 52
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 23.809523809523
$ws.Range("L17").Value = 1.960784313725
$ws.Range("M17").Value = 173.684210526316
$ws.Range("N17").Value = 8.333333333333
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = -5.882352941176
$ws.Range("L18").Value = -11.111111111111
$ws.Range("M18").Value = -51.515151515151
$ws.Range("N18").Value = -92.380952380952
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -41.666666666666
$ws.Range("I19").Value = 61
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = -23.75
$ws.Range("L19").Value = -34.408602150537
$ws.Range("M19").Value = 41.860465116279
$ws.Range("N19").Value = -27.380952380952
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 43
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -27.118644067796
$ws.Range("L20").Value = -10.416666666666
$ws.Range("M20").Value = -2.272727272727
$ws.Range("N20").Value = -92.153284671532
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 124
$ws.Range("H21").Value = -10.483870967741
$ws.Range("I21").Value = 202
$ws.Range("J21").Value = 232
$ws.Range("K21").Value = -12.931034482758
$ws.Range("L21").Value = -16.872427983539
$ws.Range("M21").Value = 4.123711340206
$ws.Range("N21").Value = -80.798479087452
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 3
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = 32.5
$ws.Range("I24").Value = 190
$ws.Range("J24").Value = 161
$ws.Range("K24").Value = 18.012422360248
$ws.Range("L24").Value = 12.426035502958
$ws.Range("M24").Value = 123.529411764706
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 61
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 48.780487804878
$ws.Range("I25").Value = 105
$ws.Range("J25").Value = 71
$ws.Range("K25").Value = 47.887323943662
$ws.Range("L25").Value = 43.835616438356
$ws.Range("C26").Value = 10
$ws.Range("E26").Value = 66.666666666666
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 12.903225806451
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 67
$ws.Range("K26").Value = 13.432835820895
$ws.Range("L26").Value = 5.555555555555
$ws.Range("M26").Value = 11.764705882352
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 33.333333333333
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -37.5
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -22.222222222222
